# Update the "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to match the refreshed scrape output.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F9").Value = 124
$wsExhibit.Range("F11").Value = 1034
$wsExhibit.Range("F15").Value = 6248
$wsExhibit.Range("F18").Value = 136
$wsExhibit.Range("F20").Value = 15035
$wsExhibit.Range("F21").Value = 1493
$wsExhibit.Range("F22").Value = 261
$wsExhibit.Range("F25").Value = 10945
$wsExhibit.Range("F27").Value = 4274

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 124
$wsAll.Range("F13").Value = 1034
$wsAll.Range("F18").Value = 6248
$wsAll.Range("F21").Value = 136
$wsAll.Range("F23").Value = 15035
$wsAll.Range("F24").Value = 1493
$wsAll.Range("F25").Value = 261
$wsAll.Range("F28").Value = 10945
$wsAll.Range("F30").Value = 4274
